$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.71"
$ws.Range("E2").Value = "'-0.65%"
$ws.Range("D3").Value = "'40.49"
$ws.Range("E3").Value = "'0.88%"
$ws.Range("D4").Value = "'5.004"
$ws.Range("E4").Value = "'-0.53%"
$ws.Range("D5").Value = "'0.07382"
$ws.Range("E5").Value = "'0.02%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.573"
$ws.Range("E6").Value = "'0.63%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9237"
$ws.Range("E7").Value = "'0.13%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.378"
$ws.Range("E8").Value = "'-0.88%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1189"
$ws.Range("E9").Value = "'-0.05%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1810"
$ws.Range("E10").Value = "'1.99%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04400"
$ws.Range("E11").Value = "'5.58%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08809"
$ws.Range("E12").Value = "'0.75%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("E13").Value = "'0.04%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001268"
$ws.Range("E14").Value = "'-0.64%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005809"
$ws.Range("E15").Value = "'0.52%"
$ws.Range("B16").Value = "HotbitToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D16").Value = "'0.003781"
$ws.Range("E16").Value = "'-3.16%"
$ws.Range("D17").Value = "'3.342"
$ws.Range("E17").Value = "'-1.18%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.293"
$ws.Range("E18").Value = "'-0.27%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3318"
$ws.Range("E19").Value = "'0.66%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'7.938"
$ws.Range("E20").Value = "'5.14%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1390"
$ws.Range("E21").Value = "'3.39%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2960"
$ws.Range("E22").Value = "'5.40%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.03918"
$ws.Range("E23").Value = "'2.90%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001262"
$ws.Range("E24").Value = "'-1.78%"
$ws.Range("E25").Value = "'-3.22%"
$ws.Range("E26").Value = "'-0.14%"
$ws.Range("D38").Value = "'0.02331"
$ws.Range("E38").Value = "'0.79%"
$ws.Range("D39").Value = "'0.05066"
$ws.Range("E39").Value = "'0.60%"
$ws.Range("D40").Value = "'0.005836"
$ws.Range("E40").Value = "'32.14%"
$ws.Range("D41").Value = "'0.007806"
$ws.Range("E41").Value = "'0.88%"
$ws.Range("D42").Value = "'0.1290"
$ws.Range("E42").Value = "'1.13%"
$ws.Range("E43").Value = "'-0.28%"
$ws.Range("D44").Value = "'0.008032"
$ws.Range("E44").Value = "'15.05%"
$ws.Range("D45").Value = "'0.2914"
$ws.Range("E45").Value = "'-8.83%"
$ws.Range("D46").Value = "'0.00006208"
$ws.Range("E46").Value = "'-3.98%"
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("D48").Value = "'0.04627"
$ws.Range("E48").Value = "'-81.62%"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("E51").Value = "'-0.12%"
